$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.291.18"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "1.839.52"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9980"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6698"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9989"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07427"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2963"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07720"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.034"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.790.78"
$ws.Range("E13").Value = "  -3.11%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6802"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.212"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.36%  "
$ws.Range("D17").Value = "28.995.00"
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008261"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "230.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.57"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9977"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.275"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9988"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "160.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.737"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1415"
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.513"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.220"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.103"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.203"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05342"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.872"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.79%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7538"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.140"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.678"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").Value = "1.331.79"
$ws.Range("E37").Value = "  +2.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01807"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.731"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9236"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.994"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.83%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("B43").Value = "XinFinNetwork"
$ws.Range("C43").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.08247"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +19.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "103.14"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.10%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000123"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.64%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5171"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "64.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.06%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.769"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "1.931.06"
$ws.Range("E49").Value = "  -2.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.296"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05949"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.26%  "
